# Auto-generated update of cached market-price values in the Chocobo Profits workbook.
# Each worksheet (one per crafting class) is an Excel Table (Table_<CLASS>) of Leve
# data; columns H-N hold cached FFXIV marketboard price/profit figures refreshed by
# the scheduled runner. We only update the cached numeric values - no structural change.

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H74").Value = 4766816.5
$ws.Range("I74").Value = 5559633
$ws.Range("J74").Value = 9919.666999999999
$ws.Range("K74").Value = 5559633
$ws.Range("L74").Value = 9919.666999999999
$ws.Range("M74").Value = -5558697
$ws.Range("N74").Value = -11791.667
$ws.Range("H77").Value = 4766816.5
$ws.Range("I77").Value = 5559633
$ws.Range("J77").Value = 9919.666999999999
$ws.Range("K77").Value = 27798165
$ws.Range("L77").Value = 49598.335
$ws.Range("M77").Value = -27793485
$ws.Range("N77").Value = -58958.335
$ws.Range("H98").Value = 3022.8462
$ws.Range("I98").Value = 1931.2727
$ws.Range("J98").Value = 9026.5
$ws.Range("K98").Value = 1931.2727
$ws.Range("L98").Value = 9026.5
$ws.Range("M98").Value = -433.2727
$ws.Range("N98").Value = -12022.5
$ws.Range("H113").Value = 9162
$ws.Range("I113").Value = 3933.3333
$ws.Range("K113").Value = 3933.3333
$ws.Range("M113").Value = -679.3332999999998
$ws.Range("H122").Value = 3022.8462
$ws.Range("I122").Value = 1931.2727
$ws.Range("J122").Value = 9026.5
$ws.Range("K122").Value = 5793.8181
$ws.Range("L122").Value = 27079.5
$ws.Range("M122").Value = -3343.8181
$ws.Range("N122").Value = -31979.5
$ws.Range("H124").Value = 42580
$ws.Range("J124").Value = 42580
$ws.Range("L124").Value = 42580
$ws.Range("N124").Value = -52400
$ws.Range("H126").Value = 42780
$ws.Range("J126").Value = 42780
$ws.Range("L126").Value = 42780
$ws.Range("N126").Value = -52660

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 665.5143
$ws.Range("I2").Value = 618.125
$ws.Range("J2").Value = 768.9091
$ws.Range("K2").Value = 618.125
$ws.Range("L2").Value = 768.9091
$ws.Range("M2").Value = -505.125
$ws.Range("N2").Value = -994.9091
$ws.Range("H32").Value = 9850.271000000001
$ws.Range("I32").Value = 7549.25
$ws.Range("J32").Value = 11428.114
$ws.Range("K32").Value = 7549.25
$ws.Range("L32").Value = 11428.114
$ws.Range("M32").Value = -7262.25
$ws.Range("N32").Value = -12002.114
$ws.Range("H45").Value = 1399.7142
$ws.Range("I45").Value = 866.7273
$ws.Range("J45").Value = 1986
$ws.Range("K45").Value = 866.7273
$ws.Range("L45").Value = 1986
$ws.Range("M45").Value = -489.7273
$ws.Range("N45").Value = -2740
$ws.Range("H63").Value = 5543664.5
$ws.Range("I63").Value = 10656924
$ws.Range("J63").Value = 4300
$ws.Range("K63").Value = 10656924
$ws.Range("L63").Value = 4300
$ws.Range("M63").Value = -10656238
$ws.Range("N63").Value = -5672
$ws.Range("H66").Value = 5543664.5
$ws.Range("I66").Value = 10656924
$ws.Range("J66").Value = 4300
$ws.Range("K66").Value = 53284620
$ws.Range("L66").Value = 21500
$ws.Range("M66").Value = -53281188
$ws.Range("N66").Value = -28364
$ws.Range("H102").Value = 1535
$ws.Range("I102").Value = 1705
$ws.Range("J102").Value = 1450
$ws.Range("K102").Value = 1705
$ws.Range("L102").Value = 1450
$ws.Range("M102").Value = -83
$ws.Range("N102").Value = -4694
$ws.Range("H116").Value = 665.5143
$ws.Range("I116").Value = 618.125
$ws.Range("J116").Value = 768.9091
$ws.Range("K116").Value = 618.125
$ws.Range("L116").Value = 768.9091
$ws.Range("M116").Value = 1675.875
$ws.Range("N116").Value = -5356.9091
$ws.Range("H122").Value = 3028.25
$ws.Range("I122").Value = 1770.25
$ws.Range("J122").Value = 4286.25
$ws.Range("K122").Value = 5310.75
$ws.Range("L122").Value = 12858.75
$ws.Range("M122").Value = -2860.75
$ws.Range("N122").Value = -17758.75
$ws.Range("H124").Value = 21750
$ws.Range("J124").Value = 21750
$ws.Range("L124").Value = 21750
$ws.Range("N124").Value = -31570
$ws.Range("H125").Value = 42248.332
$ws.Range("J125").Value = 42248.332
$ws.Range("L125").Value = 42248.332
$ws.Range("N125").Value = -52088.332
$ws.Range("H132").Value = 3059.8157
$ws.Range("I132").Value = 2233.1482
$ws.Range("J132").Value = 5088.909
$ws.Range("K132").Value = 6699.444600000001
$ws.Range("L132").Value = 15266.727
$ws.Range("M132").Value = -4169.444600000001
$ws.Range("N132").Value = -20326.727
$ws.Range("H137").Value = 51548
$ws.Range("J137").Value = 51548
$ws.Range("L137").Value = 51548
$ws.Range("N137").Value = -61748

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 665.5143
$ws.Range("I3").Value = 618.125
$ws.Range("J3").Value = 768.9091
$ws.Range("K3").Value = 618.125
$ws.Range("L3").Value = 768.9091
$ws.Range("M3").Value = -504.125
$ws.Range("N3").Value = -996.9091
$ws.Range("H105").Value = 2703
$ws.Range("I105").Value = 2712.6667
$ws.Range("J105").Value = 2596.6667
$ws.Range("K105").Value = 2712.6667
$ws.Range("L105").Value = 2596.6667
$ws.Range("M105").Value = -965.6667000000002
$ws.Range("N105").Value = -6090.6667

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H122").Value = 2380.8635
$ws.Range("I122").Value = 1730.4736
$ws.Range("J122").Value = 6500
$ws.Range("K122").Value = 5191.4208
$ws.Range("L122").Value = 19500
$ws.Range("M122").Value = -2741.4208
$ws.Range("N122").Value = -24400

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H12").Value = 2244720.8
$ws.Range("I12").Value = 10101140
$ws.Range("J12").Value = 29.5
$ws.Range("K12").Value = 30303420
$ws.Range("L12").Value = 88.5
$ws.Range("M12").Value = -30303247
$ws.Range("N12").Value = -434.5

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H123").Value = 10966.611
$ws.Range("J123").Value = 10966.611
$ws.Range("L123").Value = 10966.611
$ws.Range("N123").Value = -15866.611
$ws.Range("H124").Value = 42113.332
$ws.Range("J124").Value = 42113.332
$ws.Range("L124").Value = 42113.332
$ws.Range("N124").Value = -51933.332
$ws.Range("H137").Value = 72652.89
$ws.Range("J137").Value = 72652.89
$ws.Range("L137").Value = 72652.89
$ws.Range("N137").Value = -82852.89

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H81").Value = 60349.668
$ws.Range("J81").Value = 60349.668
$ws.Range("L81").Value = 60349.668
$ws.Range("N81").Value = -62345.668
$ws.Range("H84").Value = 60349.668
$ws.Range("J84").Value = 60349.668
$ws.Range("L84").Value = 181049.004
$ws.Range("N84").Value = -191033.004
$ws.Range("H122").Value = 3592.6667
$ws.Range("I122").Value = 2655.0527
$ws.Range("J122").Value = 12500
$ws.Range("K122").Value = 7965.158100000001
$ws.Range("L122").Value = 37500
$ws.Range("M122").Value = -5515.158100000001
$ws.Range("N122").Value = -42400
$ws.Range("H132").Value = 3803.2258
$ws.Range("I132").Value = 1530.6072
$ws.Range("K132").Value = 4591.821599999999
$ws.Range("M132").Value = -2061.821599999999

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H113").Value = 211.41667
$ws.Range("I113").Value = 212.45454
$ws.Range("J113").Value = 200
$ws.Range("K113").Value = 637.3636200000001
$ws.Range("L113").Value = 600
$ws.Range("M113").Value = 1532.63638
$ws.Range("N113").Value = -4940
$ws.Range("H122").Value = 5106.7646
$ws.Range("I122").Value = 3927.5715
$ws.Range("J122").Value = 5932.2
$ws.Range("K122").Value = 11782.7145
$ws.Range("L122").Value = 17796.6
$ws.Range("M122").Value = -9332.7145
$ws.Range("N122").Value = -22696.6

